# Update the "2024" sheet: a new SMS/notification entry was logged for
# September at the top of the (reverse-chronological) log block that lives
# in columns R:S starting at row 31. This pushes all the existing log rows
# in that block (and the two blocks below it, in columns P:Q and column A)
# down by one row - i.e. a plain row insert at row 31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above the current row 31; everything at/after row 31
# (through the old row 72) shifts down to rows 32-73, which also grows the
# sheet's used range from A1:Y72 to A1:Y73 - matching the diff exactly.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new September log entry.
$ws.Cells.Item(31, 18).Value = "bal axisbank"
$ws.Cells.Item(31, 19).Value = "2024-09-05 16:26:55"
